{"js": "// Locate the paragraphs we need by their text content so the edit is\n// resilient to exact index assumptions.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet emptyAfterHeading = null; // empty paragraph right after \"Projektidee, Beschreibung:\"\nlet dataVisParagraph = null;  // the \"data visulaization\" bullet paragraph\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text === \"Projektidee, Beschreibung:\" && i + 1 < paragraphs.items.length) {\n    emptyAfterHeading = paragraphs.items[i + 1];\n  }\n  if (text === \"data visulaization\") {\n    dataVisParagraph = paragraphs.items[i];\n  }\n}\n\n// 1) Remove the empty placeholder paragraph right after the\n//    \"Projektidee, Beschreibung:\" heading.\nif (emptyAfterHeading) {\n  emptyAfterHeading.delete();\n  await context.sync();\n}\n\n// 2) Add a new empty bullet paragraph (same sub-level as the existing\n//    \"statistical graphics\" / \"data visulaization\" bullets) right after\n//    \"data visulaization\". Splitting off of that paragraph naturally\n//    inherits its list formatting (ilvl=1, numId=1) and \"Normal\" style.\nif (dataVisParagraph) {\n  dataVisParagraph.insertParagraph(\"\", \"After\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Step 1: remove the empty placeholder paragraph right after the\n#     \"Projektidee, Beschreibung:\" heading. ---------------------------------\n$emptyAfterHeading = $null\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.TrimEnd(\"`r\", \"`n\")\n    if ($t -eq \"Projektidee, Beschreibung:\" -and ($i + 1) -le $count) {\n        $emptyAfterHeading = $d.Paragraphs.Item($i + 1)\n        break\n    }\n}\nif ($emptyAfterHeading -ne $null) {\n    $emptyAfterHeading.Range.Delete()\n}\n\n# --- Step 2: add a new empty bullet paragraph (same sub-level as the\n#     existing \"statistical graphics\" / \"data visulaization\" bullets) right\n#     after \"data visulaization\". Splitting off of that paragraph naturally\n#     inherits its list formatting (ilvl=1, numId=1) and \"Normal\" style.\n#     Re-query the paragraph collection fresh (post step 1) so we don't use\n#     a stale reference captured before the document shifted. -------------\n$dataVisParagraph = $null\n$count2 = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count2; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.TrimEnd(\"`r\", \"`n\")\n    if ($t -eq \"data visulaization\") {\n        $dataVisParagraph = $p\n        break\n    }\n}\nif ($dataVisParagraph -ne $null) {\n    $dataVisParagraph.Range.InsertParagraphAfter()\n}\n"}
